# Lotomania data update: append 6 newly-drawn lottery contests (2882-2887)
# to the LOTOMANIA worksheet, as rows 376-381 (columns A:U).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is one contest: contest number followed by the 20 drawn numbers.
$newContests = @(
    @(2882, 1, 3, 6, 9, 14, 15, 16, 19, 23, 31, 32, 41, 44, 63, 64, 66, 71, 77, 91, 97),
    @(2883, 3, 5, 13, 16, 21, 28, 33, 35, 40, 43, 47, 53, 61, 63, 71, 73, 76, 77, 81, 89),
    @(2884, 0, 3, 8, 14, 23, 30, 32, 37, 38, 41, 45, 55, 58, 60, 71, 74, 77, 78, 80, 91),
    @(2885, 0, 10, 13, 14, 18, 27, 32, 42, 43, 44, 49, 52, 55, 60, 61, 69, 73, 76, 81, 82),
    @(2886, 1, 3, 5, 12, 13, 28, 35, 48, 49, 50, 72, 74, 75, 76, 77, 80, 81, 82, 90, 92),
    @(2887, 0, 2, 4, 5, 13, 14, 20, 30, 36, 49, 69, 76, 77, 78, 81, 85, 88, 90, 96, 97)
)

$startRow = 376
for ($i = 0; $i -lt $newContests.Count; $i++) {
    $rowValues = $newContests[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowValues[$c]
    }
}

$lastRow = $startRow + $newContests.Count - 1

# Mirror the author's final selection over the newly pasted block.
[void]$ws.Range("B" + $startRow + ":U" + $lastRow).Select()

Write-Host "Added contests $startRow to $lastRow"
